$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must remain TEXT (not a number,
# which would lose precision/format). Assigning a numeric-looking string
# directly via .Value converts it to a number and/or forces a NumberFormat
# change that creates a new style. Instead, compute the text in a scratch
# cell via a formula (naturally yields a text result) and paste back only
# the value, which preserves B3's original style (s="8").
$scratch = $ws.Range("ZZ1000")
$scratch.Formula = '="2570314725427075"'
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 12.04.2024"

$ws.Range("B6").Value = "14.04."
$ws.Range("C6").Value = "15.04."
$ws.Range("D6").Value = "PAYPAL PADFKX"
$ws.Range("E6").Value = "17,38-"

$ws.Range("B7").Value = "17.04."
$ws.Range("C7").Value = "18.04."
$ws.Range("D7").Value = "KARTENZ./17.04 ALDI SUED RO"
$ws.Range("E7").Value = "103,58-"

$ws.Range("B8").Value = "21.04."
$ws.Range("C8").Value = "22.04."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-99970195"
$ws.Range("E8").Value = "57,10-"

$ws.Range("B9").Value = "23.04."
$ws.Range("C9").Value = "24.04."
$ws.Range("D9").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E9").Value = "47,46-"

$ws.Range("B10").Value = "27.04."
$ws.Range("C10").Value = "28.04."
$ws.Range("D10").Value = "AMAZON.DE MKTPLC EU BXSKJM"
$ws.Range("E10").Value = "60,96-"

$ws.Range("B11").Value = "29.04."
$ws.Range("C11").Value = "30.04."
$ws.Range("D11").Value = "ZALANDO MKTPLC EU GGRJBY"
$ws.Range("E11").Value = "204,13-"

$ws.Range("D12").Value = "KONTOSTAND AM 01.05.2024"
$ws.Range("E12").Value = "490,61-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 11.05.2024"
